{"js": "// Fix regenerate prompt in batch mode\n// Fill in the first still-empty trailing row of the \"Ollama Model\" results\n// table with the qwen3-next:80b-cloud entry.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Target table is the second table in the document (the Ollama model\n// results table: \"Ollama Model\" | \"Status\" | \"Notes\").\nconst resultsTable = tables.items[1];\n\nconst rows = resultsTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Find the first row (after the header) whose cells are all still empty.\nlet targetRow = null;\nfor (let i = 1; i < rows.items.length; i++) {\n  const cells = rows.items[i].cells;\n  cells.load(\"items/body/text\");\n  await context.sync();\n\n  const isEmpty = cells.items.every((c) => c.body.text.trim() === \"\");\n  if (isEmpty) {\n    targetRow = rows.items[i];\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error(\"No empty row found to fill in.\");\n}\n\nconst cells = targetRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\n// Column 1: Ollama model name.\ncells.items[0].body.insertText(\"qwen3-next:80b-cloud\", \"Replace\");\n\n// Column 2: Status \u2014 \"1st try worked (some errors on prompt re-gen)\"\n// with \"st\" rendered as a superscript, matching the other rows' style\n// (e.g. \"2nd\", \"3rd\").\nconst statusBody = cells.items[1].body;\nstatusBody.insertText(\"1\", \"Replace\");\nconst superscriptRun = statusBody.insertText(\"st\", \"End\");\nsuperscriptRun.font.superscript = true;\nstatusBody.insertText(\" try worked\", \"End\");\nstatusBody.insertText(\" (some errors on prompt re-gen)\", \"End\");\n\n// Column 3: Notes \u2014 timing info.\ncells.items[2].body.insertText(\"~2min for 24 prompts\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Fix regenerate prompt in batch mode\n# Fill in the first still-empty trailing row of the \"Ollama Model\" results\n# table with the qwen3-next:80b-cloud entry.\n\n$d = $word.ActiveDocument\n\n# Locate the results table (header: \"Ollama Model\" | \"Status\" | \"Notes\").\n$targetTable = $null\nfor ($ti = 1; $ti -le $d.Tables.Count; $ti++) {\n    $tbl = $d.Tables.Item($ti)\n    $headerText = $tbl.Rows.Item(1).Cells.Item(1).Range.Text\n    if ($headerText.StartsWith(\"Ollama Model\")) {\n        $targetTable = $tbl\n    }\n}\n\n# Find the first row (after the header) whose cells are all still empty.\n# An \"empty\" cell's Range.Text is just the paragraph mark + end-of-cell\n# mark (length 2), since it has no real text runs yet.\n$targetRowIndex = 0\nfor ($ri = 2; $ri -le $targetTable.Rows.Count; $ri++) {\n    $row = $targetTable.Rows.Item($ri)\n    $c1 = $row.Cells.Item(1).Range.Text\n    $c2 = $row.Cells.Item(2).Range.Text\n    $c3 = $row.Cells.Item(3).Range.Text\n    if ($c1.Length -le 2 -and $c2.Length -le 2 -and $c3.Length -le 2) {\n        $targetRowIndex = $ri\n        break\n    }\n}\n\n$row = $targetTable.Rows.Item($targetRowIndex)\n\n# Column 1: Ollama model name.\n$row.Cells.Item(1).Range.Text = \"qwen3-next:80b-cloud\"\n\n# Column 2: Status \u2014 \"1st try worked (some errors on prompt re-gen)\"\n# with \"st\" rendered as a superscript, matching the other rows' style\n# (e.g. \"2nd\", \"3rd\").\n$cell2 = $row.Cells.Item(2)\n$cell2.Range.Text = \"1st try worked\"\n$cell2Start = $cell2.Range.Start\n$d.Range($cell2Start + 1, $cell2Start + 3).Font.Superscript = $true\n\n$cell2End = $cell2.Range.End - 1\n$d.Range($cell2End, $cell2End).InsertBefore(\" (some errors on prompt re-gen)\")\n\n# Column 3: Notes \u2014 timing info.\n$row.Cells.Item(3).Range.Text = \"~2min for 24 prompts\"\n"}
